# Swap the START/END water node id values (columns B and C) for the rows
# where the edge direction was corrected, then scroll the sheet view down
# to show the bottom of the data (to compare the metrics / plot distribution).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToSwap = @(2,4,5,7,8,14,15,17,18,24,25,26,27,28,29,32,33,34,44,45,48,53,58,59,60,61,62,64,65,69,71)

foreach ($r in $rowsToSwap) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value2 = $cVal
    $ws.Cells.Item($r, 3).Value2 = $bVal
}

# Scroll the view so row 52 is the top-left visible row and select B72,
# matching where the author ended up after reviewing the data.
$excel.ActiveWindow.ScrollRow = 52
$ws.Range("B72").Select()
